# Update data to 02-06, predict 02-7
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("chinancoron2020")

# New row 22
$ws.Range("A22").Value = 43866
$ws.Range("B22").Value = 28060
$ws.Range("C22").Value = 3697
$ws.Range("D22").Value = 24702
$ws.Range("E22").Value = 564

# New row 23
$ws.Range("A23").Value = 43866
$ws.Range("B23").Value = 31211
$ws.Range("C23").Value = 3151
$ws.Range("D23").Value = 26359
$ws.Range("E23").Value = 637

# Match the existing date style (numFmtId 14) used by column A, reusing the
# same style index instead of minting a new custom number format.
$ws.Range("A21").Copy()
$ws.Range("A22:A23").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Update the selected cell to mirror the author's final cursor position
$ws.Range("G17").Select()
